$d = $word.ActiveDocument

# The last paragraph currently ends: "...defined on its own line " followed
# by the (hidden) _GoBack bookmark. We need to:
#   1. split off a new, empty, plain paragraph after it
#   2. split off a second new, plain paragraph holding "New test line"
#   3. move the _GoBack bookmark (collapsed) to the very end of that text

# _GoBack is Word's special "last edit" bookmark; it's hidden from the
# Bookmarks collection's Count/enumeration but is still reachable by name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- split 1: new empty paragraph after "...own line " ---
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara1.Range.ListFormat.RemoveNumbers()
$newPara1.Style = "Normal"

# --- split 2: new paragraph that will hold "New test line" ---
$endRange2 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange2.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2.Range.ListFormat.RemoveNumbers()
$newPara2.Style = "Normal"

# A collapsed Range sitting exactly one char before the document's very
# last position mis-resolves when handed to Bookmarks.Add in this runtime,
# so insert the text together with one placeholder trailing character,
# anchor the bookmark there (where it is no longer the very-last position),
# then delete the placeholder; the collapsed bookmark stays put.
$endRange3 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange3.InsertAfter("New test line#")

$bookmarkPos = $d.Content.End - 2
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$placeholder.Delete()

Write-Output "done"
